$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update columns C:G for rows 2-25
$cg = New-Object 'object[,]' 24,5
$cg[0,0] = 13.06314679895733
$cg[0,1] = 3.621259098539409
$cg[0,2] = 13.09130027399132
$cg[0,3] = 76.22404913125168
$cg[0,4] = 3.883766919405143
$cg[1,0] = 13.06387408571281
$cg[1,1] = 3.566338061518021
$cg[1,2] = 13.13103634300009
$cg[1,3] = 75.74530099706223
$cg[1,4] = 3.889386187930529
$cg[2,0] = 13.06745735867809
$cg[2,1] = 3.531829946880816
$cg[2,2] = 13.15774147947367
$cg[2,3] = 75.4622248886846
$cg[2,4] = 3.893007898387516
$cg[3,0] = 13.06970530709803
$cg[3,1] = 3.517573406734413
$cg[3,2] = 13.16920425839003
$cg[3,3] = 75.34964758521664
$cg[3,4] = 3.89452709535095
$cg[4,0] = 13.07012611694119
$cg[4,1] = 3.515194493927039
$cg[4,2] = 13.17114269072114
$cg[4,3] = 75.33112319659654
$cg[4,4] = 3.894781979342082
$cg[5,0] = 13.06748448740192
$cg[5,1] = 3.53163845982154
$cg[5,2] = 13.15789372108436
$cg[5,3] = 75.4606953279562
$cg[5,4] = 3.893028211131221
$cg[6,0] = 13.06274604471843
$cg[6,1] = 3.60248763727524
$cg[6,2] = 13.10452250575586
$cg[6,3] = 76.05674660626244
$cg[6,4] = 3.885668981960449
$cg[7,0] = 13.07837255258504
$cg[7,1] = 3.735033711727678
$cg[7,2] = 13.01816326036257
$cg[7,3] = 77.30980883611255
$cg[7,4] = 3.872588558846667
$cg[8,0] = 13.1050633227541
$cg[8,1] = 3.828330535261768
$cg[8,2] = 12.96586813488159
$cg[8,3] = 78.27870304998308
$cg[8,4] = 3.86378874871367
$cg[9,0] = 13.12050375988957
$cg[9,1] = 3.869848694437065
$cg[9,2] = 12.94449871272491
$cg[9,3] = 78.72922663195736
$cg[9,4] = 3.85995860180809
$cg[10,0] = 13.12682380226268
$cg[10,1] = 3.885434669619085
$cg[10,2] = 12.93675459206538
$cg[10,3] = 78.90116509074612
$cg[10,4] = 3.858532868959269
$cg[11,0] = 13.12544165227967
$cg[11,1] = 3.882084053199524
$cg[11,2] = 12.93840694834511
$cg[11,3] = 78.86407683575189
$cg[11,4] = 3.85883883254319
$cg[12,0] = 13.12101424105609
$cg[12,1] = 3.871133707137484
$cg[12,2] = 12.9438546243027
$cg[12,3] = 78.7433457677533
$cg[12,4] = 3.859840812858287
$cg[13,0] = 13.11836388841023
$cg[13,1] = 3.864408490179201
$cg[13,2] = 12.94723680609975
$cg[13,3] = 78.66956623337711
$cg[13,4] = 3.860457760136144
$cg[14,0] = 13.10412053383895
$cg[14,1] = 3.825598365066702
$cg[14,2] = 12.96731337465447
$cg[14,3] = 78.2494509484946
$cg[14,4] = 3.864042519985539
$cg[15,0] = 13.09622680876096
$cg[15,1] = 3.80155079219905
$cg[15,2] = 12.98024952012209
$cg[15,3] = 77.99417713676957
$cg[15,4] = 3.866285798967959
$cg[16,0] = 13.09199713713487
$cg[16,1] = 3.787632573848172
$cg[16,2] = 12.98791781865298
$cg[16,3] = 77.84827424944537
$cg[16,4] = 3.867592363167665
$cg[17,0] = 13.09061841836241
$cg[17,1] = 3.782905318632811
$cg[17,2] = 12.99055329119515
$cg[17,3] = 77.79903478960317
$cg[17,4] = 3.868037547367574
$cg[18,0] = 13.09703497210311
$cg[18,1] = 3.804119694017928
$cg[18,2] = 12.97884887176218
$cg[18,3] = 78.02125635032495
$cg[18,4] = 3.866045313688475
$cg[19,0] = 13.12230185192189
$cg[19,1] = 3.87435380732353
$cg[19,2] = 12.94224506448859
$cg[19,3] = 78.7787717488348
$cg[19,4] = 3.859545839317395
$cg[20,0] = 13.14157181969496
$cg[20,1] = 3.91946142002311
$cg[20,2] = 12.92035087707505
$cg[20,3] = 79.2815987153362
$cg[20,4] = 3.855441701960193
$cg[21,0] = 13.131035375886
$cg[21,1] = 3.895460387672516
$cg[21,2] = 12.93185059600038
$cg[21,3] = 79.01254542755618
$cg[21,4] = 3.857619082206933
$cg[22,0] = 13.09666864034471
$cg[22,1] = 3.80295858276552
$cg[22,2] = 12.9794813849511
$cg[22,3] = 78.00901116192482
$cg[22,4] = 3.866153984511685
$cg[23,0] = 13.07147380743333
$cg[23,1] = 3.69988042661816
$cg[23,2] = 13.0395668315569
$cg[23,3] = 76.96210747347912
$cg[23,4] = 3.875983891998546
$ws.Range("C2:G25").Value = $cg

# Update columns J:L for rows 2-25
$jl = New-Object 'object[,]' 24,3
$jl[0,0] = 13.10892412913287
$jl[0,1] = 31.4627909636547
$jl[0,2] = 9.132339507478726
$jl[1,0] = 13.12548267589665
$jl[1,1] = 31.39887296689587
$jl[1,2] = 9.159109560404513
$jl[2,0] = 13.13731426532832
$jl[2,1] = 31.36984443130212
$jl[2,2] = 9.17645863598349
$jl[3,0] = 13.142553537913
$jl[3,1] = 31.3605897277458
$jl[3,2] = 9.183758915549078
$jl[4,0] = 13.14344872644542
$jl[4,1] = 31.35920859793114
$jl[4,2] = 9.184985066385979
$jl[5,0] = 13.13738323334032
$jl[5,1] = 31.36970918926045
$jl[5,2] = 9.176556155901169
$jl[6,0] = 13.11428757892843
$jl[6,1] = 31.43863550336361
$jl[6,2] = 9.141381119529832
$jl[7,0] = 13.08223912938144
$jl[7,1] = 31.65450765010221
$jl[7,2] = 9.079594285834306
$jl[8,0] = 13.06681314207211
$jl[8,1] = 31.8616469017022
$jl[8,2] = 9.038519312312836
$jl[9,0] = 13.06156828158572
$jl[9,1] = 31.9662211122077
$jl[9,2] = 9.020757518118891
$jl[10,0] = 13.05983777947897
$jl[10,1] = 32.00728735274014
$jl[10,2] = 9.014163298518332
$jl[11,0] = 13.06019909437904
$jl[11,1] = 31.9983781584561
$jl[11,2] = 9.015577636364721
$jl[12,0] = 13.0614207852206
$jl[12,1] = 31.96957040216479
$jl[12,2] = 9.020212372234344
$jl[13,0] = 13.06220241551962
$jl[13,1] = 31.9521150909693
$jl[13,2] = 9.023068414591478
$jl[14,0] = 13.06719163640393
$jl[14,1] = 31.85501913648792
$jl[14,2] = 9.03969858335306
$jl[15,0] = 13.07070687132839
$jl[15,1] = 31.79808922028671
$jl[15,2] = 9.050136450739169
$jl[16,0] = 13.072895536058
$jl[16,1] = 31.76631985553207
$jl[16,2] = 9.056227012126056
$jl[17,0] = 13.07366520464232
$jl[17,1] = 31.75573139979746
$jl[17,2] = 9.058304140260432
$jl[18,0] = 13.07031540057378
$jl[18,1] = 31.80404872291805
$jl[18,2] = 9.049016327487253
$jl[19,0] = 13.06105500210822
$jl[19,1] = 31.97799232780406
$jl[19,2] = 9.018847469279814
$jl[20,0] = 13.05649301692967
$jl[20,1] = 32.10020932115278
$jl[20,2] = 8.999898064841522
$jl[21,0] = 13.05879124691539
$jl[21,1] = 32.03420659346747
$jl[21,2] = 9.009941805205305
$jl[22,0] = 13.07049186218049
$jl[22,1] = 31.80135143874587
$jl[22,2] = 9.049522455762059
$jl[23,0] = 13.08948641723871
$jl[23,1] = 31.58752571443546
$jl[23,2] = 9.095546025670108
$ws.Range("J2:L25").Value = $jl
